# Slide 42 ("Code Generation for CompoundStmt") - Content Placeholder 2:
# The grammar-rule line currently reads:
#     compoundStmt = "{" statement "}" .
# It needs to become (and be split into three separate runs, matching the
# author's edit):
#     compoundStmt = "{" statements "}" .
#   run A: ' = '
#   run B: '"{" statements '
#   run C: '"}" .'

$p = $ppt.ActivePresentation

# Find the slide / shape / paragraph that holds the grammar rule, rather than
# hard-coding indices, so the script is resilient to small renumbering.
$targetSlide = $null
$targetShape = $null
$targetParaIndex = -1

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }
        $tr = $tf.TextRange

        $paraCount = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            if ($para.Text -like '*compoundStmt = "{" statement "}" .*') {
                $targetSlide = $slide
                $targetShape = $shape
                $targetParaIndex = $pi
                break
            }
        }
        if ($targetParaIndex -ge 1) { break }
    }
    if ($targetParaIndex -ge 1) { break }
}

$tr = $targetShape.TextFrame.TextRange
$para = $tr.Paragraphs($targetParaIndex, 1)

$oldSnippet = ' = "{" statement "}" .'
$newSnippet = ' = "{" statements "}" .'

$relIdx = $para.Text.IndexOf($oldSnippet)
$absStart = $para.Start + $relIdx

# First, update the text in place (still a single run at this point).
$snippetRange = $tr.Characters($absStart, $oldSnippet.Length)
$snippetRange.Text = $newSnippet

# Now split that single run into three runs with the exact text the author
# ended up with. Re-assigning .Text on each sub-range forces PowerPoint to
# break the run apart at that boundary while keeping the existing run
# formatting (font, size, etc.) on each of the resulting pieces.
$seg1 = ' = '
$seg2 = '"{" statements '
$seg3 = '"}" .'

$run1 = $tr.Characters($absStart, $seg1.Length)
$run1.Text = $seg1

$run2 = $tr.Characters($absStart + $seg1.Length, $seg2.Length)
$run2.Text = $seg2

$run3 = $tr.Characters($absStart + $seg1.Length + $seg2.Length, $seg3.Length)
$run3.Text = $seg3
